$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 73

# Write the date label as a formula first so Excel's auto date-recognition
# doesn't kick in (it would otherwise turn "01-07-2021" into a date serial
# and mint a new number-format style). Then convert it in place to a
# static value via copy / paste-special so it lands as a plain shared
# string, matching the rest of column A.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Formula = '="01-07-2021"'
$dateCell.Copy() | Out-Null
$dateCell.PasteSpecial(-4163) | Out-Null

$ws.Cells.Item($row, 2).Value = 60.91
$ws.Cells.Item($row, 3).Value = 50.44
$ws.Cells.Item($row, 4).Value = 62.05
$ws.Cells.Item($row, 5).Value = 61.7
$ws.Cells.Item($row, 6).Value = 59.56
$ws.Cells.Item($row, 7).Value = 59.06
